$d = $word.ActiveDocument

function Get-ParaIndexLike($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# --- 1. Replace the "Typical example: ..." paragraph with the new explanatory text ---
$idx = Get-ParaIndexLike("Typical example:*")
$p = $d.Paragraphs.Item($idx)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = 'Command to execute starts at the first arg not starting with "-", with {} replaced by one line, surrounded by quotes if arg contains spaces. '

# Insert the two extra explanation paragraphs right after it (new paragraphs come back
# empty, so fill them in via InsertAfter rather than a Range.Text= assignment).
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($idx + 1)
$newp.Range.InsertAfter("If there is no {} in command, suppose there is one as last argument (append arg to command).")

$p2 = $d.Paragraphs.Item($idx + 1)
$p2.Range.InsertParagraphAfter()
$newp2 = $d.Paragraphs.Item($idx + 2)
$newp2.Range.InsertAfter("If standard input is empty, do not execute command (silently unless verbose mode)")

# --- 2. After the following empty paragraph, insert "Example" and the new rfind/rgrep/rxargs example paragraph ---
$emptyIdx = $idx + 3   # the already-existing empty paragraph right after the 3 paragraphs above
$emptyPara = $d.Paragraphs.Item($emptyIdx)
$emptyPara.Range.InsertParagraphAfter()
$examplep = $d.Paragraphs.Item($emptyIdx + 1)
$examplep.Range.InsertAfter("Example")

$examplep = $d.Paragraphs.Item($emptyIdx + 1)
$examplep.Range.InsertParagraphAfter()
$cmdp = $d.Paragraphs.Item($emptyIdx + 2)
$cmdp.Range.InsertAfter("rfind C:\development\Git*\**\*.rs -exec rgrep -l pomme {} | rxargs notepad++ {}")

# Insert a fresh blank paragraph before the "Command line options:" paragraph (cleaner than
# appending one, since it leaves that paragraph's own runs untouched).
$optIdx = Get-ParaIndexLike("Command line options:*")
$optPara = $d.Paragraphs.Item($optIdx)
$startOfOpt = $d.Range($optPara.Range.Start, $optPara.Range.Start)
$startOfOpt.InsertParagraphBefore()

# --- 3. Rename "Command line options:" header to "rxargs options:" ---
$optIdx = Get-ParaIndexLike("Command line options:*")
$optPara = $d.Paragraphs.Item($optIdx)
$headerRange = $d.Range($optPara.Range.Start, $optPara.Range.Start + "Command line options:".Length)
$headerRange.Text = "rxargs options:"

# --- 4. Extend the -@ file option description ---
$atIdx = Get-ParaIndexLike("-@ file, instead of reading standard input, read @file*")
$atPara = $d.Paragraphs.Item($atIdx)
$atEnd = $d.Range($atPara.Range.End - 1, $atPara.Range.End - 1)
$atEnd.InsertAfter(". Similar to redirecting standard input, but using textautodecode crate to support most text file variants")

# --- 5. Extend the -1 option description ---
$oneIdx = Get-ParaIndexLike("-1, instead of executing one command per file*")
$onePara = $d.Paragraphs.Item($oneIdx)
$oneEnd = $d.Range($onePara.Range.End - 1, $onePara.Range.End - 1)
$oneEnd.InsertAfter(" (or multiple commands ensuring than arguments size does not exceed 7800 UTF-16 chars)")

# --- 6. Append a new final paragraph about rfind / ActionPrint / stdout flush ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Range.InsertAfter("Possibly update rfind to call after each ActionPrint: io::stdout().flush().unwrap(); to avoid being blocked by output buffering")
